$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row-level changes -----------------------------------------------
# Old row 11 ("MH100, MH101, ... / M3 / CubeSat_Backplane:M3_PTH_PAD") is removed entirely.
$ws.Rows(11).Delete()

# A new row is inserted at row 10 for the J101 burn-wire connector.
$ws.Rows(10).Insert()

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "J101"
$ws.Range("F10").Value = "DF13-4P-1.25V(76)"
$ws.Range("E10").Value = "CubeSat_Backplane:DF13-4P-1.25(75)"
$ws.Range("D10").Value = "BURN_WIRES"

# --- New footprint labels for the solar-panel connectors (rows 6-8) --
# Typed with a leading apostrophe so they are stored as text with the
# quote-prefix flag set, matching how Excel marks this kind of entry.
$ws.Range("D6").Value = "'X+ SOLAR PANEL"
$ws.Range("D7").Value = "'X- SOLAR PANEL"
$ws.Range("D8").Value = "'Y+ SOLAR PANEL"

# --- DNI note on row 19 (SW100, SW101 / RESET) ------------------------
$ws.Range("G19").Value = "DNI"

# --- Selection moves to D4 --------------------------------------------
$null = $ws.Range("D4").Select()

Write-Host "edits applied"
